$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("NodeShapes (classes)")
$ws3 = $wb.Worksheets.Item("PropertyShapes (properties)")
$ws2.Activate()
$ws2.Range("B13").Select()
$ws3.Activate()
